# Update workbook/worksheet for the "Through 2021-10-23" data refresh
# (adds one more day of data, 2021-10-31 run covering through 10-23)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2021-10-23"

# Update the row label for October
$ws.Range("A12").Value = "October (through 10-23)"

# Row 12 (October partial month) updates
$ws.Range("C12").Value = 18
$ws.Range("D12").Value = 0.1
$ws.Range("F12").Value = 36
$ws.Range("G12").Value = 0.0769
$ws.Range("I12").Value = 43
$ws.Range("J12").Value = 0.14
$ws.Range("L12").Value = 49
$ws.Range("M12").Value = 0.0577
$ws.Range("O12").Value = 35
$ws.Range("P12").Value = 0.1026
$ws.Range("R12").Value = 111
$ws.Range("U12").Value = 147

# Row 13 (Total) updates
$ws.Range("C13").Value = 214
$ws.Range("D13").Value = 0.1301
$ws.Range("F13").Value = 419
$ws.Range("G13").Value = 0.1047
$ws.Range("I13").Value = 620
$ws.Range("J13").Value = 0.0842
$ws.Range("L13").Value = 536
$ws.Range("M13").Value = 0.1067
$ws.Range("O13").Value = 414
$ws.Range("P13").Value = 0.102
$ws.Range("R13").Value = 959
$ws.Range("S13").Value = 0.0524
$ws.Range("U13").Value = 1312
$ws.Range("V13").Value = 0.0588
